# Refresh the cryptos worksheet with the latest scraped Price (column D)
# and Volume(1h) (column E) figures, matching a GitHub Actions scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as plain text (e.g. "1.827.17"). Some refreshed values
# look like ordinary decimal numbers (e.g. "241.28"), and Excel would silently
# reinterpret a bare assignment as a number (dropping trailing zeros / exact
# text form). Forcing the number format to Text ("@") before the assignment,
# then restoring the "Normal" style afterwards, keeps the cell as plain text
# while leaving its original (default) style/formatting untouched.
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "29.051.52"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.825.91"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue $ws.Range("D5") "241.28"
$ws.Range("E5").Value = "  +0.04%  "
Set-TextValue $ws.Range("D6") "0.6352"
$ws.Range("E6").Value = "  -4.59%  "
Set-TextValue $ws.Range("D7") "1.001"
$ws.Range("E7").Value = "  +0.03%  "
Set-TextValue $ws.Range("D8") "44.80"
$ws.Range("E8").Value = "  +6.74%  "
Set-TextValue $ws.Range("D9") "0.2932"
$ws.Range("E9").Value = "  +0.15%  "
Set-TextValue $ws.Range("D10") "0.07333"
$ws.Range("E10").Value = "  -0.58%  "
Set-TextValue $ws.Range("D11") "22.77"
$ws.Range("E11").Value = "  +0.45%  "
Set-TextValue $ws.Range("D12") "0.07669"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").Value = "1.827.17"
$ws.Range("E13").Value = "  -0.49%  "
Set-TextValue $ws.Range("D14") "4.982"
$ws.Range("E14").Value = "  -0.05%  "
Set-TextValue $ws.Range("D15") "0.6627"
$ws.Range("E15").Value = "  -0.96%  "
Set-TextValue $ws.Range("D16") "81.93"
$ws.Range("E16").Value = "  -1.25%  "
Set-TextValue $ws.Range("D17") "0.000008665"
$ws.Range("E17").Value = "  +4.77%  "
Set-TextValue $ws.Range("D18") "6.028"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").Value = "29.055.15"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").Value = "2.074.31"
$ws.Range("E20").Value = "  +0.21%  "
Set-TextValue $ws.Range("D21") "224.57"
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("E22").Value = "  -0.64%  "
Set-TextValue $ws.Range("D23") "1.001"
Set-TextValue $ws.Range("D24") "7.111"
$ws.Range("E24").Value = "  -0.27%  "
Set-TextValue $ws.Range("D25") "1.001"
$ws.Range("E25").Value = "  +0.05%  "
Set-TextValue $ws.Range("D26") "158.47"
$ws.Range("E26").Value = "  -1.44%  "
Set-TextValue $ws.Range("D27") "8.459"
$ws.Range("E27").Value = "  -1.92%  "
Set-TextValue $ws.Range("D28") "0.1368"
$ws.Range("E28").Value = "  -1.74%  "
Set-TextValue $ws.Range("D29") "17.88"
$ws.Range("E29").Value = "  -0.48%  "
Set-TextValue $ws.Range("D30") "1.504"
$ws.Range("E30").Value = "  -0.44%  "
Set-TextValue $ws.Range("D31") "4.086"
$ws.Range("E31").Value = "  -0.66%  "
Set-TextValue $ws.Range("D32") "4.023"
$ws.Range("E32").Value = "  -0.28%  "
Set-TextValue $ws.Range("D34") "0.05286"
$ws.Range("E34").Value = "  -0.35%  "
Set-TextValue $ws.Range("D35") "1.835"
$ws.Range("E35").Value = "  -1.90%  "
Set-TextValue $ws.Range("D36") "0.7359"
$ws.Range("E36").Value = "  -2.29%  "
$ws.Range("E37").Value = "  +2.03%  "
Set-TextValue $ws.Range("D38") "2.650"
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("D39").Value = "1.295.66"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("E40").Value = "  -0.63%  "
Set-TextValue $ws.Range("D41") "2.743"
$ws.Range("E41").Value = "  +0.83%  "
Set-TextValue $ws.Range("D42") "6.309"
$ws.Range("E42").Value = "  +5.84%  "
Set-TextValue $ws.Range("D43") "0.8996"
$ws.Range("E43").Value = "  -2.18%  "
Set-TextValue $ws.Range("D44") "1.0000"
$ws.Range("E44").Value = "  -0.73%  "
Set-TextValue $ws.Range("D45") "102.43"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "1.974.04"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("E47").Value = "  -0.49%  "
Set-TextValue $ws.Range("D48") "63.97"
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("E50").Value = "  -2.26%  "
Set-TextValue $ws.Range("D51") "0.07263"
$ws.Range("E51").Value = "  -17.03%  "
